$d = $word.ActiveDocument

# En dash (U+2013) and right single quotation mark (U+2019) used in the
# source text below.
$enDash = [char]0x2013
$rsquo  = [char]0x2019

# Grab the Word bullet-list template (same one the "Bullets" button in the
# ribbon applies) so the three new paragraphs share one numbering
# definition, just like selecting them together and clicking "Bullets".
$bulletGallery  = $word.ListGalleries.Item(1)
$bulletTemplate = $bulletGallery.ListTemplates.Item(1)

# --- Insert the three new list paragraphs after the existing paragraph ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$p1 = $d.Paragraphs.Last.Range
$p1.Text = "Region Salary Tuition Income.csv " + $enDash + " contains colleges" + $rsquo + " programs and costs of program and the income level for that program."
$p1.Style = "List Paragraph"

$p1end = $d.Paragraphs.Last.Range
$p1end.InsertParagraphAfter()

$p2 = $d.Paragraphs.Last.Range
$p2.Text = "Seshu_Miriyala_Final.Rmd" + " " + $enDash + " R script analyzing the college data."
$p2.Style = "List Paragraph"

$p2end = $d.Paragraphs.Last.Range
$p2end.InsertParagraphAfter()

$p3 = $d.Paragraphs.Last.Range
$p3.Text = "Seshu_Miriyala_Final.html" + " " + $enDash + " Report generated post analysis."
$p3.Style = "List Paragraph"

# Apply one shared bulleted-list numbering definition across all three
# new paragraphs (mirrors selecting them and clicking the Bullets button).
$firstNewPara = $d.Paragraphs.Item($d.Paragraphs.Count - 2)
$newParasRange = $d.Range($firstNewPara.Range.Start, $p3.End)
$newParasRange.ListFormat.ApplyListTemplate($bulletTemplate)
